$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.090.36"
$ws.Range("E2").Value = "  -3.92%  "

$ws.Range("D3").Value = "1.961.79"
$ws.Range("E3").Value = "  -6.49%  "

$ws.Range("E4").Value = "  +0.82%  "

$ws.Range("D5").Value = "'327.33"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -4.57%  "

$ws.Range("E6").Value = "  +0.81%  "

$ws.Range("D7").Value = "'0.4990"
$ws.Range("D7").ClearFormats()

$ws.Range("D8").Value = "'0.4204"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -4.26%  "

$ws.Range("D9").Value = "'54.11"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -1.63%  "

$ws.Range("D10").Value = "'0.09067"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -3.33%  "

$ws.Range("D11").Value = "'1.096"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -6.77%  "

$ws.Range("D12").Value = "'22.96"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -7.27%  "

$ws.Range("D13").Value = "1.929.15"
$ws.Range("E13").Value = "  -8.28%  "

$ws.Range("D14").Value = "'7.853"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -8.35%  "

$ws.Range("D15").Value = "'6.416"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -6.65%  "

$ws.Range("D16").Value = "'1.012"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.97%  "

$ws.Range("D17").Value = "'0.00001099"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -5.02%  "

$ws.Range("D18").Value = "'90.96"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -10.32%  "

$ws.Range("D19").Value = "'0.06651"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.25%  "

$ws.Range("D20").Value = "'19.17"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -9.26%  "

$ws.Range("D21").Value = "'1.010"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.78%  "

$ws.Range("D22").Value = "'5.936"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -7.65%  "

$ws.Range("D23").Value = "29.125.45"
$ws.Range("E23").Value = "  -3.86%  "

$ws.Range("D24").Value = "'11.94"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.97%  "

$ws.Range("D25").Value = "'2.293"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -1.22%  "

$ws.Range("D26").Value = "'156.62"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -3.44%  "

$ws.Range("D27").Value = "'20.58"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.68%  "

$ws.Range("D28").Value = "'6.198"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -11.38%  "

$ws.Range("D29").Value = "'2.253"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -10.69%  "

$ws.Range("D30").Value = "'126.70"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.17%  "

$ws.Range("D31").Value = "'1.038"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -8.84%  "

$ws.Range("D32").Value = "'0.09837"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -6.60%  "

$ws.Range("D33").Value = "'1.523"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -9.13%  "

$ws.Range("D34").Value = "'5.804"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -7.35%  "

$ws.Range("D35").Value = "'3.712"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -4.12%  "

$ws.Range("D36").Value = "'0.02426"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -7.87%  "

$ws.Range("D37").Value = "'8.986"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -11.16%  "

$ws.Range("D38").Value = "'0.06332"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -6.43%  "

$ws.Range("D39").Value = "'1.285"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -4.50%  "

$ws.Range("D40").Value = "'0.6430"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -7.70%  "

$ws.Range("D41").Value = "'11.44"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -9.78%  "

$ws.Range("D42").Value = "'0.1996"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -9.88%  "

$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("D44").Value = "'0.6194"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -8.87%  "

$ws.Range("D45").Value = "'13.39"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -6.69%  "

$ws.Range("D46").Value = "'2.171"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -6.74%  "

$ws.Range("D47").Value = "'1.287"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.65%  "

$ws.Range("D48").Value = "'3.470"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.51%  "

$ws.Range("E49").Value = "  -4.30%  "

$ws.Range("D50").Value = "'0.06858"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -5.93%  "

$ws.Range("D51").Value = "'1.104"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -8.75%  "
